$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$invoiceLong = "INVOICE IN THE NAME OF NOHAR CHAND LOHA BHANDAR , UMESH R SHARMA, KATHLAL ROAD, KANBHA, AHMEDABAD, Gujarat, 382430 NOHAR CHAND LOHA BHANDAR , UMESH R SHARMA, KATHLAL ROAD, KANBHA, AHMEDABAD, Gujarat, 382430. KKF - KANKARIYA W.RLY"
$invoiceShort = "INVOICE IN THE NAME OF"

# Rows 2-4: clear Consignee Name (H) and set Destination (I) to long invoice text
foreach ($r in 2..4) {
    $ws.Cells.Item($r, 8).Value = ""
    $ws.Cells.Item($r, 9).Value = $invoiceLong
}

# Rows 5-6: clear Consignee Name (H) and set Destination (I) to short invoice text
foreach ($r in 5..6) {
    $ws.Cells.Item($r, 8).Value = ""
    $ws.Cells.Item($r, 9).Value = $invoiceShort
}
